$d = $word.ActiveDocument

# --- 1. "Blockchain Development / Intern - United Health Group - Optum" ---
# Merge " Intern - United Health Group - " and "Optum" runs (removes proofErr split)
# without touching the separate "Blockchain Development" run.
$r = $d.Content
$r.Find.Execute(" Intern " + [char]8211 + " United Health Group " + [char]8211 + " Optum", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rr = $d.Range($r.Start, $r.End)
$rr.Text = " Intern " + [char]8211 + " United Health Group " + [char]8211 + " Optum"

# --- 2. "Made use of Hyperledger Fabric and Composer ..." ---
# Merge "Hyperledger" + " Fabric and Composer to develop permissioned Blockchain solutions"
$r = $d.Content
$r.Find.Execute("Hyperledger Fabric and Composer to develop permissioned Blockchain solutions", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rr = $d.Range($r.Start, $r.End)
$rr.Text = "Hyperledger Fabric and Composer to develop permissioned Blockchain solutions"

# --- 3. "Uses Golang, Javascript, D3, Nodejs, Bash, Docker, Python, and Git." ---
# Whole paragraph -> merges into one run, then re-split into the two runs the diff wants:
# "Uses Golang, Javascrip" | "t, D3, Nodejs, Bash, Docker, Python, and Git."
$r = $d.Content
$r.Find.Execute("Uses Golang, Javascript, D3, Nodejs, Bash, Docker, Python, and Git.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$paraStart = $r.Start
$rr = $d.Range($r.Start, $r.End)
$rr.Text = "Uses Golang, Javascript, D3, Nodejs, Bash, Docker, Python, and Git."
$splitPos = $paraStart + ("Uses Golang, Javascrip").Length
$bm = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bm)

# --- 4. "Created Composer Blockchain deployment pipeline with an interactive Blockchain visualizer" ---
# Capitalize lowercase "blockchain" -> "Blockchain" (the semantic fix), producing 3 runs:
# "Created ... interactive " | "B" | "lockchain visualizer"
$r = $d.Content
$r.Find.Execute("interactive blockchain visualizer", $true, $false, $false, $false, $false, $true, 1, $false, "interactive Blockchain visualizer", 2)

$r2 = $d.Content
$r2.Find.Execute("Created Composer Blockchain deployment pipeline with an interactive Blockchain visualizer", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$paraStart2 = $r2.Start
$splitPos1 = $paraStart2 + ("Created Composer Blockchain deployment pipeline with an interactive ").Length
$bm = $d.Range($splitPos1, $splitPos1)
$d.Bookmarks.Add("_GoBack", $bm)
$splitPos2 = $paraStart2 + ("Created Composer Blockchain deployment pipeline with an interactive B").Length
$bm = $d.Range($splitPos2, $splitPos2)
$d.Bookmarks.Add("_GoBack", $bm)

# --- 5. "portable graphics libraries to support arbitrary data visualizations." ---
# Lowercase "Graphics Libraries" -> "graphics libraries", then re-split runs to match
# original neighbouring run boundaries ("portable " / "graphics libraries" / " to support...")
# and relocate the _GoBack bookmark to sit mid-word ("graphics l" | "ibraries"),
# matching where the _GoBack bookmark ends up after this was the very last edit made.
$r = $d.Content
$r.Find.Execute("Graphics Libraries", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$glStart = $r.Start
$rr = $d.Range($r.Start, $r.End)
$rr.Text = "graphics libraries"

$r3 = $d.Content
$r3.Find.Execute("portable graphics libraries to support", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$portableEnd = $r3.Start + ("portable ").Length
$bm = $d.Range($portableEnd, $portableEnd)
$d.Bookmarks.Add("_GoBack", $bm)

$glEnd = $glStart + ("graphics libraries").Length
$bm = $d.Range($glEnd, $glEnd)
$d.Bookmarks.Add("_GoBack", $bm)

# Final _GoBack position: between "graphics l" and "ibraries" (this was the last edit)
$midPos = $glStart + ("graphics l").Length
$bm = $d.Range($midPos, $midPos)
$d.Bookmarks.Add("_GoBack", $bm)

# --- 6. "BadgerBlockchain - Cofounder" ---
$r = $d.Content
$r.Find.Execute("BadgerBlockchain " + [char]8211 + " Cofounder", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rr = $d.Range($r.Start, $r.End)
$rr.Text = "BadgerBlockchain " + [char]8211 + " Cofounder"

# --- 7. "Made use of Python / Flask, Solidity, Bash, and Javascript for implementation" ---
$r = $d.Content
$r.Find.Execute("Made use of Python / Flask, Solidity, Bash, and Javascript for implementation", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rr = $d.Range($r.Start, $r.End)
$rr.Text = "Made use of Python / Flask, Solidity, Bash, and Javascript for implementation"

# --- 8. "Python, Bash Scripting, Nodejs, and the Twitter API ... Mbps. " ---
# Merge these runs, leaving the preceding separate "Uses " run untouched.
$r = $d.Content
$r.Find.Execute("Python, Bash Scripting, Nodejs, and the Twitter API to monitor internet speed and tweet at Charter whenever internet service drops below advertised download Mbps. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rr = $d.Range($r.Start, $r.End)
$rr.Text = "Python, Bash Scripting, Nodejs, and the Twitter API to monitor internet speed and tweet at Charter whenever internet service drops below advertised download Mbps. "
